$d = $word.ActiveDocument

# --- Work Experience: QA bullet rewording ---
# "Led QA teams to perform automated and manual testing on various projects."
# becomes "Performed automated and manual testing on various projects within many QA teams."
$d.Content.Find.Execute(
    "Led QA teams to perform automated and manual testing on various projects.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Performed automated and manual testing on various projects within many QA teams.",
    2) | Out-Null

# --- Personal Projects: repurpose "Personal Website - Test Automation (Cypress)" heading ---
# into the new "Stock Tracker - Video Game Consoles" project heading. Replace the bold run
# and the plain (non-bold) run separately so the bold/non-bold run split survives intact.
$enDash = [char]8211
$d.Content.Find.Execute(
    "Personal Website " + $enDash + " Test Automation ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Stock Tracker - Video Game Consoles ",
    2) | Out-Null

$d.Content.Find.Execute(
    "(Cypress)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "(Java, Selenium WebDriver with Cucumber)",
    2) | Out-Null

# --- Personal Projects: bullet describing the stock tracker project ---
$oldBullet1 = "Created a suite of automated test cases using Cypress that verifies personal website and projects are stable, functional, and contains correct information"
$newBullet1 = "Created a stock tracker that automatically verifies the current stock of newly released consoles. Major retailers it tracks include Amazon and Best Buy"
$d.Content.Find.Execute(
    $oldBullet1,
    $true, $false, $false, $false, $false, $true, 1, $false,
    $newBullet1,
    2) | Out-Null

# --- Personal Projects: "Weather Application (JavaScript, HTML, CSS)" heading becomes the
# "Test Automation Suite for Personal Website (Cypress)" heading (old Cypress project, moved
# down). Again replace the bold run and the plain run separately to keep the bold/non-bold
# run split intact (and to drop the tab character that used to precede the parenthetical).
$tab = [char]9
$d.Content.Find.Execute(
    "Weather Application" + $tab + " ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Test Automation Suite for Personal Website ",
    2) | Out-Null

$d.Content.Find.Execute(
    "(JavaScript, HTML, CSS)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "(Cypress)",
    2) | Out-Null

# --- Personal Projects: the two Weather Application bullets collapse into a single bullet
# (the second one is removed outright, the first one is reworded). ---
$weatherBullet2 = "Takes current data from OpenWeatherMap" + [char]8217 + "s API and displays any city" + [char]8217 + "s forecast"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $weatherBullet2) {
        $p.Range.Delete() | Out-Null
        break
    }
}

$oldBullet2 = "Created a functional responsive weather application using AJAX with JSON"
$newBullet2 = "Created a suite of automated test cases that verifies personal website and projects are stable, functional, and contains expected information"
$d.Content.Find.Execute(
    $oldBullet2,
    $true, $false, $false, $false, $false, $true, 1, $false,
    $newBullet2,
    2) | Out-Null
